# Backlog.xlsx - "Estructura de producto en tabla de articulos"
#
# Applies the page-layout / view tweaks captured in the target edit:
#   - Sheet1: turn on "Fit to page" (sheetPr/pageSetUpPr fitToPage) and
#     scale print output to 84%
#   - Sheet1: narrower left/right print margins (0.7in -> 0.25in)
#   - Sheet1: move the active selection from C76 to A74
#   - Workbook window position/size bookkeeping (xWindow/yWindow/size)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Page setup: fit-to-page + 84% scale ---------------------------------
# Turning on "Fit to" first (so sheetPr/pageSetUpPr@fitToPage is recorded),
# then applying the 84% scale on top of it.
$ws.PageSetup.Zoom = 84
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1

# --- Page margins: 0.7in -> 0.25in on left/right -------------------------
$ws.PageSetup.LeftMargin = $excel.InchesToPoints(0.25)
$ws.PageSetup.RightMargin = $excel.InchesToPoints(0.25)

# --- Move the selected / active cell from C76 to A74 ---------------------
$ws.Range("A74").Select()

# --- Workbook window geometry (cosmetic UI bookkeeping) ------------------
$win = $wb.Windows.Item(1)
$win.Left = 240
$win.Top = 645
$win.Width = 14805
$win.Height = 7470
